$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Annisa Putri Restu's row (row 8): Repayment_collections, Repayment_amount,
# and Pending Amount Recovery change; Pending Amount stays the same value.
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = "595,048.00"
$ws.Range("G8").Value = "0.31"

# Rename the worksheet/tab to reflect the new upload revision.
$ws.Name = "repayment_20250915_20250915 (3)"
